# perubahan 11 juni 2025 penyesuaian penamaan sheet pada template afi dan msi
#
# Renames the (single) worksheet from "Home Credit Indonesia" to "TigerSnus"
# (this also repoints the _xlnm._FilterDatabase defined name, since it is
# expressed relative to the sheet name) and resets the stored cursor
# selection on that sheet back to A1.

$wb = $excel.ActiveWorkbook

# Target the workbook's first/only sheet explicitly rather than relying on
# whichever one happens to be "active".
$ws = $wb.Worksheets.Item(1)

# 1) Rename "Home Credit Indonesia" -> "TigerSnus".
$ws.Name = "TigerSnus"

# 2) Reset the saved selection on that sheet from H17 back to A1.
$ws.Range("A1").Select()
